# "Generate Report for Handoff"
# The cbdfbb95-... file has moved from "In Translation" to "Ready for handoff",
# with refreshed handoff timestamps and a translation-type change (ht -> mt).
# Update the Overview sheet and each per-locale sheet (zh-cn, de-de) accordingly.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the cbdfbb95-... file ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-23 10:14:27"

# --- zh-cn sheet: row 3 is the cbdfbb95-... file ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-23 10:14:23"

# --- de-de sheet: row 3 is the cbdfbb95-... file ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-23 10:14:27"
